# Generate Report for Handoff
# Update the localization-status report: the two source files that were
# previously handed back are replaced by a new pair of files that are now
# ready for handoff (not yet handed back), with fresh handoff timestamps
# and a fresh target (.xlf) file per language. The now-unused
# "Latest Target File" / "Latest Handback File" columns are cleared since
# no handback has happened yet for the new files.

$wb = $excel.ActiveWorkbook

# Old source files were "8c158229-aaf3-4380-9354-f2773b6c84aa.md" and
# "e39c1c27-11e0-450a-82b9-3da1447aa221.md"; they're superseded below.
$newMd1 = "1e2b5805-6f60-4125-a897-ef3151d8ab4c.md"
$newMd2 = "ffffe73ca54b-970e-4d2a-a723-8512024a563e.md"

$newStatus = "Ready for handoff"

$newXlfZhCn = "1e2b5805-6f60-4125-a897-ef3151d8ab4c.056ac7fcbe6e14b6529a7349561b36bd236bafa1.zh-cn.xlf"
$newXlfDeDe = "1e2b5805-6f60-4125-a897-ef3151d8ab4c.056ac7fcbe6e14b6529a7349561b36bd236bafa1.de-de.xlf"

$newHandoffDtZhCn = "2016-03-23 05:12:15"
$newHandoffDtDeDe = "2016-03-23 05:12:19"

$newHandbackDt = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Drop the now-stale "Latest Target File" (F) / "Latest Handback File" (G)
# columns for both data rows - nothing has been handed back yet.
$wsZh.Range("F2:G3").Clear()

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Range("D2").Value = $newXlfZhCn
$wsZh.Range("D3").Value = $newXlfZhCn

$wsZh.Range("E2").Value = $newHandoffDtZhCn
$wsZh.Range("E3").Value = $newHandoffDtZhCn

$wsZh.Range("H2").Value = $newHandbackDt
$wsZh.Range("H3").Value = $newHandbackDt

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$newMd1", "", "", $newMd1)
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newXlfZhCn", "", "", $newXlfZhCn)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$newMd2", "", "", $newMd2)
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newXlfZhCn", "", "", $newXlfZhCn)

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("F2:G3").Clear()

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Range("D2").Value = $newXlfDeDe
$wsDe.Range("D3").Value = $newXlfDeDe

$wsDe.Range("E2").Value = $newHandoffDtDeDe
$wsDe.Range("E3").Value = $newHandoffDtDeDe

$wsDe.Range("H2").Value = $newHandbackDt
$wsDe.Range("H3").Value = $newHandbackDt

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$newMd1", "", "", $newMd1)
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newXlfDeDe", "", "", $newXlfDeDe)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$newMd2", "", "", $newMd2)
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newXlfDeDe", "", "", $newXlfDeDe)

# ---------------------------------------------------------------------
# Overview sheet - File Name hyperlinks refresh to the new source files,
# the zh-cn/de-de status columns both flip to the new status, and the
# rollup "Latest Handoff Date" picks up the newest handoff timestamp
# (the de-de one).
# ---------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

$wsOv.Range("B2").Value = $newStatus
$wsOv.Range("C2").Value = $newStatus
$wsOv.Range("D2").Value = $newHandoffDtDeDe

$wsOv.Range("B3").Value = $newStatus
$wsOv.Range("C3").Value = $newStatus
$wsOv.Range("D3").Value = $newHandoffDtDeDe

$wsOv.Hyperlinks.Delete()
$wsOv.Hyperlinks.Add($wsOv.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$newMd1", "", "", $newMd1)
$wsOv.Hyperlinks.Add($wsOv.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$newMd2", "", "", $newMd2)
